# Lagt til ekstra linje mellom bus 1 og 2
# Inserts a new data row right after row 2 on sheet "BranchData", duplicating
# row 2's values/format into the new row 3, and shifting the old rows 3-10
# down to rows 4-11 (their values/styles travel with them unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 3..10 down to rows 4..11.
# Go from the bottom up so we never overwrite a row before it has been copied.
for ($r = 10; $r -ge 3; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    $src = $ws.Range("A" + $srcRow + ":E" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":E" + $dstRow)
    $src.Copy($dst)
}

# New row 3 duplicates row 2 (same line/bus data and formatting).
$ws.Range("A2:E2").Copy($ws.Range("A3:E3"))

# Match the author's final cursor position from the recorded session.
[void]$ws.Range("H8").Select()
